$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "2025-04-28 11:45:56"
$ws.Range("B11").Value = 226
